$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows for the additional "2508" period entries ---
# (pushes the old rows 22 onward, i.e. the closing signature block, down by 2)
$ws.Range("22:23").Insert()

# --- 2. Fix up the formatting of the table rows ---
# Row 21 currently still carries the old "bottom of table" (outer-border)
# styling. Grab a copy of that look first, so it can be moved onto the new
# last row (23) further down.
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
# Row 21 is no longer the last row of the table (two more rows follow it
# now), so give it the same look as the other interior rows (e.g. row 20)
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
# Row 22 (newly inserted) should also look like a normal data row
$ws.Range("B22:J22").PasteSpecial(-4122)

# --- 3. Rewrite the whole worker/period table (rows 16-23), now sorted
#        by period (2505-2508), two rows (CC + PPT) per period ---
$ws.Cells.Item(16,2).Value = "CC"
$ws.Cells.Item(16,3).Value = "1143379475"
$ws.Cells.Item(16,4).Value = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Cells.Item(16,5).Value = "2505"

$ws.Cells.Item(17,2).Value = "PPT"
$ws.Cells.Item(17,3).Value = "5614268"
$ws.Cells.Item(17,4).Value = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Cells.Item(17,5).Value = "2505"

$ws.Cells.Item(18,2).Value = "CC"
$ws.Cells.Item(18,3).Value = "1143379475"
$ws.Cells.Item(18,4).Value = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Cells.Item(18,5).Value = "2506"

$ws.Cells.Item(19,2).Value = "PPT"
$ws.Cells.Item(19,3).Value = "5614268"
$ws.Cells.Item(19,4).Value = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Cells.Item(19,5).Value = "2506"

$ws.Cells.Item(20,2).Value = "CC"
$ws.Cells.Item(20,3).Value = "1143379475"
$ws.Cells.Item(20,4).Value = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Cells.Item(20,5).Value = "2507"

$ws.Cells.Item(21,2).Value = "PPT"
$ws.Cells.Item(21,3).Value = "5614268"
$ws.Cells.Item(21,4).Value = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Cells.Item(21,5).Value = "2507"

$ws.Cells.Item(22,2).Value = "CC"
$ws.Cells.Item(22,3).Value = "1143379475"
$ws.Cells.Item(22,4).Value = "LILIBETH ESTHER PEREZ BABILONIA"
$ws.Cells.Item(22,5).Value = "2508"
$ws.Cells.Item(22,6).Value = 56940
$ws.Cells.Item(22,7).Value = 1423500

$ws.Cells.Item(23,2).Value = "PPT"
$ws.Cells.Item(23,3).Value = "5614268"
$ws.Cells.Item(23,4).Value = "JEFFERSON SMITH RIVERA CASTILLO"
$ws.Cells.Item(23,5).Value = "2508"
$ws.Cells.Item(23,6).Value = 56940
$ws.Cells.Item(23,7).Value = 1423500

# --- 4. Update the account-summary totals ---
# Total overdue amount (Valor Mora) now covers 4 periods instead of 3
$ws.Range("E11").Value = 455520
# Number of overdue periods (Cant. Periodos)
$ws.Range("F13").Value = 4
